# Update MSME definitions: replace literal "<br/>" separators with real line breaks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")
$nl = [char]10

$ws.Range("B24").Value = "<10 Manufacturing" + $nl + "<5 Other businesses"
$ws.Range("C24").Value = "< N$ 500,000 Manufacturing" + $nl + "< N$ 100,000 Other businesses"
$ws.Range("D24").Value = "< N$ 1,000,000 Manufacturing" + $nl + "< N$ 250,000 Other businesses"

$ws.Range("B25").Value = "> 10 Manufacturing" + $nl + "> 5 Other businesses"
$ws.Range("C25").Value = "> N$ 500,000 Manufacturing" + $nl + "> N$ 100,000 Other businesses"
$ws.Range("D25").Value = "> N$ 1,000,000 Manufacturing" + $nl + "> N$ 250,000 Other businesses"
